# Generate Report for Handoff
# Updates status text and timestamps across the Overview / zh-cn / de-de
# sheets to reflect that the handoff report has been (re)generated, and
# widens the "Status" columns that now hold the longer "Ready for handoff"
# text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$overview.Range("E2").Value2 = "Ready for handoff"
$overview.Range("F2").Value2 = "Ready for handoff"
$zhcn.Range("C2").Value2     = "Ready for handoff"
$dede.Range("C2").Value2     = "Ready for handoff"

# --- Timestamps -----------------------------------------------------------
# Latest HO Xliff Generate Date (Overview) / Latest Handoff Datetime (de-de)
$overview.Range("G2").Value2 = "2016-08-27 16:38:38"
$dede.Range("H2").Value2     = "2016-08-27 16:38:38"

# Latest Handoff Datetime (zh-cn)
$zhcn.Range("H2").Value2 = "2016-08-27 16:38:33"

# --- Widen the Status columns to fit the new, longer text -----------------
$overview.Columns.Item(5).ColumnWidth = 16.38265482584637
$overview.Columns.Item(6).ColumnWidth = 16.38265482584637
$zhcn.Columns.Item(3).ColumnWidth     = 16.38265482584637
$dede.Columns.Item(3).ColumnWidth     = 16.38265482584637
